{"js": "// Locate the bond (\"\ucc44\uad8c\") table: the second table in the document body,\n// with header row \"\ucc44\uad8c\uc774\ub984, \uc138\uc804\uc218\uc775\ub960, \uc138\ud6c4\uc218\uc775\ub960, \ub9cc\uae30\uc77c, \uc2e0\uc6a9\ub4f1\uae09, \ubc1c\ud589\uc0ac, \uc794\uc874\uae30\uac04(\uc77c)\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet bondTable = null;\nfor (const t of tables.items) {\n  t.load(\"values\");\n}\nawait context.sync();\nfor (const t of tables.items) {\n  if (t.values.length > 0 && t.values[0][0] === \"\ucc44\uad8c\uc774\ub984\") {\n    bondTable = t;\n    break;\n  }\n}\n\nif (!bondTable) {\n  throw new Error(\"Bond table not found\");\n}\n\nconst rows = bondTable.rows;\nrows.load(\"items/values\");\nawait context.sync();\n\n// Remove the row holding the \"\uba54\ub9ac\uce20\uce90\ud53c\ud0c8185-1\" bond (first column match).\nlet targetRow = null;\nfor (const r of rows.items) {\n  if (r.values[0][0] === \"\uba54\ub9ac\uce20\uce90\ud53c\ud0c8185-1\") {\n    targetRow = r;\n    break;\n  }\n}\nif (!targetRow) {\n  throw new Error(\"\uba54\ub9ac\uce20\uce90\ud53c\ud0c8185-1 row not found\");\n}\ntargetRow.delete();\nawait context.sync();\n\n// Append a new bond row at the end of the table for \ud558\ub098\uc5d0\ud504\uc564\uc544\uc774182-1.\nbondTable.addRows(\"End\", 1, [\n  [\"\ud558\ub098\uc5d0\ud504\uc564\uc544\uc774182-1\", \"5.0\", \"4.23\", \"2023.09.06\", \"A\", \"\ud55c\uad6d\ud22c\uc790\uc99d\uad8c\", \"217\"],\n]);\nawait context.sync();\n\n// Remove the \"\ubcf4\uc720\uc885\ubaa9 Report\" portfolio table (text + embedded chart image).\nconst tables2 = context.document.body.tables;\ntables2.load(\"items\");\nawait context.sync();\nfor (const t of tables2.items) {\n  t.load(\"values\");\n}\nawait context.sync();\n\nlet portfolioTable = null;\nfor (const t of tables2.items) {\n  if (\n    t.values.length > 0 &&\n    t.values[0][0] &&\n    t.values[0][0].indexOf(\"\ud55c\uad6d\uae08\uc735\uc9c0\uc8fc\") !== -1\n  ) {\n    portfolioTable = t;\n    break;\n  }\n}\nif (!portfolioTable) {\n  throw new Error(\"Portfolio table not found\");\n}\nportfolioTable.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the bond (\"\ucc44\uad8c\") table: the one whose header row starts with \"\ucc44\uad8c\uc774\ub984\".\n$bondTable = $null\nforeach ($t in $d.Tables) {\n    $hdr = $t.Cell(1, 1).Range.Text\n    if ($hdr -like \"*\ucc44\uad8c\uc774\ub984*\") {\n        $bondTable = $t\n        break\n    }\n}\nif ($bondTable -eq $null) {\n    throw \"Bond table not found\"\n}\n\n# Remove the row holding the \"\uba54\ub9ac\uce20\uce90\ud53c\ud0c8185-1\" bond (match on first column).\n$targetRow = $null\nforeach ($r in $bondTable.Rows) {\n    $cellText = $r.Cells.Item(1).Range.Text\n    if ($cellText -like \"*\uba54\ub9ac\uce20\uce90\ud53c\ud0c8185-1*\") {\n        $targetRow = $r\n        break\n    }\n}\nif ($targetRow -eq $null) {\n    throw \"\uba54\ub9ac\uce20\uce90\ud53c\ud0c8185-1 row not found\"\n}\n$targetRow.Delete()\n\n# Append a new bond row at the end of the table for \ud558\ub098\uc5d0\ud504\uc564\uc544\uc774182-1.\n$newRow = $bondTable.Rows.Add()\n$idx = $newRow.Index\n$bondTable.Cell($idx, 1).Range.Text = \"\ud558\ub098\uc5d0\ud504\uc564\uc544\uc774182-1\"\n$bondTable.Cell($idx, 2).Range.Text = \"5.0\"\n$bondTable.Cell($idx, 3).Range.Text = \"4.23\"\n$bondTable.Cell($idx, 4).Range.Text = \"2023.09.06\"\n$bondTable.Cell($idx, 5).Range.Text = \"A\"\n$bondTable.Cell($idx, 6).Range.Text = \"\ud55c\uad6d\ud22c\uc790\uc99d\uad8c\"\n$bondTable.Cell($idx, 7).Range.Text = \"217\"\n\n# Remove the \"\ubcf4\uc720\uc885\ubaa9 Report\" portfolio table (stock list + embedded chart image).\n$portfolioTable = $null\nforeach ($t in $d.Tables) {\n    $c1 = $t.Cell(1, 1).Range.Text\n    if ($c1 -like \"*\ud55c\uad6d\uae08\uc735\uc9c0\uc8fc*\") {\n        $portfolioTable = $t\n        break\n    }\n}\nif ($portfolioTable -eq $null) {\n    throw \"Portfolio table not found\"\n}\n$portfolioTable.Delete()\n"}
